$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '91.748.70', '  +1.26%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.133.61', '  +0.82%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.08%  ')
    ,@(5, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '241.32', '  -0.25%  ')
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '619.47', '  -0.63%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '1.11', '  -5.10%  ')
    ,@(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.390', '  +5.70%  ')
    ,@(9, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.999', '  -0.11%  ')
    ,@(10, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.130.14', '  +0.79%  ')
    ,@(11, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.753', '  -0.84%  ')
    ,@(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.204', '  +0.33%  ')
    ,@(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000255', '  +1.92%  ')
    ,@(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '35.26', '  +0.16%  ')
    ,@(15, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.62', '  +2.46%  ')
    ,@(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '91.371.04', '  +0.93%  ')
    ,@(17, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.715.97', '  +0.82%  ')
    ,@(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.130.63', '  +1.15%  ')
    ,@(19, 'SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '3.78', '  +0.48%  ')
    ,@(20, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '15.01', '  +4.33%  ')
    ,@(21, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.92', '  +1.66%  ')
    ,@(22, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '459.22', '  +2.23%  ')
    ,@(23, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0000202', '  -3.62%  ')
    ,@(24, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '9.25', '  +2.00%  ')
    ,@(25, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '5.93', '  +1.11%  ')
    ,@(26, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '89.12', '  -4.65%  ')
    ,@(27, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '1.53', '  +52.96%  ')
    ,@(28, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '11.79', '  -0.93%  ')
    ,@(29, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.151', '  +30.54%  ')
    ,@(30, 'WrappedeETH', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth', '3.298.28', '  +1.04%  ')
    ,@(31, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  -0.05%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.231', '  +1.42%  ')
    ,@(33, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.167', '  -4.94%  ')
    ,@(34, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '9.37', '  +2.74%  ')
    ,@(35, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.174', '  +8.94%  ')
    ,@(36, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '26.47', '  -1.89%  ')
    ,@(37, 'RenderToken', 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render', '7.47', '  -1.91%  ')
    ,@(38, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.96', '  +2.70%  ')
    ,@(39, 'MantraDAO', 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om', '3.94', '  -5.83%  ')
    ,@(40, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '492.98', '  +0.11%  ')
    ,@(41, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.32', '  +2.08%  ')
    ,@(42, 'PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.439', '  +5.31%  ')
    ,@(43, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '3.40', '  -5.09%  ')
    ,@(44, 'WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '22.15', '  +0.23%  ')
    ,@(45, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  -0.05%  ')
    ,@(46, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.708', '  +2.52%  ')
    ,@(47, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.93', '  +1.63%  ')
    ,@(48, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '156.17', '  -0.73%  ')
    ,@(49, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.36', '  +1.76%  ')
    ,@(50, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.50', '  -2.42%  ')
    ,@(51, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0326', '  +4.73%  ')
)

foreach ($item in $data) {
    $r = [int]$item[0]
    $ws.Cells.Item($r, 2).Value = [string]$item[1]
    $ws.Cells.Item($r, 3).Value = [string]$item[2]
    $ws.Cells.Item($r, 4).Value = [string]$item[3]
    $ws.Cells.Item($r, 5).Value = [string]$item[4]
}
